$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3123.75
$ws.Range("I18").Value = 1600
$ws.Range("J18").Value = 4647.5
$ws.Range("K18").Value = 1600
$ws.Range("L18").Value = 4647.5
$ws.Range("M18").Value = -1316
$ws.Range("N18").Value = -5215.5
$ws.Range("H40").Value = 1562.7333
$ws.Range("J40").Value = 2170.8572
$ws.Range("L40").Value = 2170.8572
$ws.Range("N40").Value = -2520.8572
$ws.Range("H100").Value = 2168.75
$ws.Range("I100").Value = 1775.1666
$ws.Range("K100").Value = 1775.1666
$ws.Range("M100").Value = -1234.1666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4093.8235
$ws.Range("I32").Value = 2379.6667
$ws.Range("K32").Value = 2379.6667
$ws.Range("M32").Value = -2092.6667
$ws.Range("H102").Value = 28573266
$ws.Range("I102").Value = 28573266
$ws.Range("K102").Value = 28573266
$ws.Range("M102").Value = -28571644

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3002.6667
$ws.Range("I20").Value = 3002.6667
$ws.Range("K20").Value = 3002.6667
$ws.Range("M20").Value = -2755.6667
$ws.Range("H86").Value = 4259.8335
$ws.Range("J86").Value = 1099.5
$ws.Range("L86").Value = 1099.5
$ws.Range("N86").Value = -3345.5
$ws.Range("H89").Value = 4259.8335
$ws.Range("J89").Value = 1099.5
$ws.Range("L89").Value = 5497.5
$ws.Range("N89").Value = -16729.5
$ws.Range("H99").Value = 5496776.5
$ws.Range("I99").Value = 6411706
$ws.Range("J99").Value = 7200
$ws.Range("K99").Value = 6411706
$ws.Range("L99").Value = 7200
$ws.Range("M99").Value = -6410208
$ws.Range("N99").Value = -10196
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H134").Value = 2581.8125
$ws.Range("I134").Value = 2581.8125
$ws.Range("K134").Value = 7745.4375
$ws.Range("M134").Value = -5210.4375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2106.25
$ws.Range("J31").Value = 1974.75
$ws.Range("L31").Value = 1974.75
$ws.Range("N31").Value = -2564.75
$ws.Range("H34").Value = 2106.25
$ws.Range("J34").Value = 1974.75
$ws.Range("L34").Value = 1974.75
$ws.Range("N34").Value = -2378.75
$ws.Range("H50").Value = 8398.1875
$ws.Range("J50").Value = 9092
$ws.Range("L50").Value = 9092
$ws.Range("N50").Value = -10342
$ws.Range("H105").Value = 1453.2858
$ws.Range("J105").Value = 1944.25
$ws.Range("L105").Value = 1944.25
$ws.Range("N105").Value = -5438.25
$ws.Range("H112").Value = 49999
$ws.Range("J112").Value = 49999
$ws.Range("L112").Value = 49999
$ws.Range("N112").Value = -52953

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1328.2142
$ws.Range("J12").Value = 2145.625
$ws.Range("L12").Value = 6436.875
$ws.Range("N12").Value = -6782.875
$ws.Range("H15").Value = 328.2857
$ws.Range("I15").Value = 328.2857
$ws.Range("K15").Value = 984.8571000000001
$ws.Range("M15").Value = -844.8571000000001
$ws.Range("H17").Value = 492.25
$ws.Range("I17").Value = 229.5
$ws.Range("J17").Value = 755
$ws.Range("K17").Value = 688.5
$ws.Range("L17").Value = 2265
$ws.Range("M17").Value = -519.5
$ws.Range("N17").Value = -2603
$ws.Range("H34").Value = 5449.7144
$ws.Range("J34").Value = 7399.8
$ws.Range("L34").Value = 22199.4
$ws.Range("N34").Value = -22367.4
$ws.Range("H39").Value = 9000
$ws.Range("J39").Value = 9000
$ws.Range("L39").Value = 27000
$ws.Range("N39").Value = -27588
$ws.Range("H55").Value = 9000
$ws.Range("J55").Value = 9000
$ws.Range("L55").Value = 27000
$ws.Range("N55").Value = -27354
$ws.Range("H62").Value = 7520.75
$ws.Range("J62").Value = 6494.3335
$ws.Range("L62").Value = 19483.0005
$ws.Range("N62").Value = -20855.0005
$ws.Range("H65").Value = 7520.75
$ws.Range("J65").Value = 6494.3335
$ws.Range("L65").Value = 58449.0015
$ws.Range("N65").Value = -65313.0015
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H123").Value = 579
$ws.Range("I123").Value = 579
$ws.Range("K123").Value = 1737
$ws.Range("M123").Value = 713
$ws.Range("H132").Value = 1988.7778
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 1999.875
$ws.Range("K132").Value = 17100
$ws.Range("L132").Value = 17998.875
$ws.Range("M132").Value = -14570
$ws.Range("N132").Value = -23058.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2655.4
$ws.Range("J80").Value = 2661.6667
$ws.Range("L80").Value = 2661.6667
$ws.Range("N80").Value = -4657.6667
$ws.Range("H83").Value = 2655.4
$ws.Range("J83").Value = 2661.6667
$ws.Range("L83").Value = 13308.3335
$ws.Range("N83").Value = -23292.3335
$ws.Range("H97").Value = 533
$ws.Range("I97").Value = 440.1
$ws.Range("K97").Value = 440.1
$ws.Range("M97").Value = 55.89999999999998
$ws.Range("H126").Value = 2325.5
$ws.Range("I126").Value = 1991
$ws.Range("K126").Value = 5973
$ws.Range("M126").Value = -3503

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2450.476
$ws.Range("I7").Value = 2536.1875
$ws.Range("K7").Value = 2536.1875
$ws.Range("M7").Value = -2424.1875
$ws.Range("H22").Value = 2393.625
$ws.Range("I22").Value = 1874.8334
$ws.Range("K22").Value = 1874.8334
$ws.Range("M22").Value = -1579.8334
$ws.Range("H27").Value = 2393.625
$ws.Range("I27").Value = 1874.8334
$ws.Range("K27").Value = 1874.8334
$ws.Range("M27").Value = -1767.8334
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H126").Value = 2450.476
$ws.Range("I126").Value = 2536.1875
$ws.Range("K126").Value = 7608.5625
$ws.Range("M126").Value = -5138.5625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1258.8
$ws.Range("I122").Value = 1323.625
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 3970.875
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -1520.875
$ws.Range("N122").Value = -7898.5
